# Refresh the crypto price table (columns B-E, rows 2-51) with the
# latest snapshot pulled by the scheduled GitHub Actions job.
#
# Only the cells whose values actually changed between snapshots are
# touched here, mirroring the upstream commit. Two rows (41/42)
# additionally swap which coin occupies them (MultiversX <-> LidoDAOToken)
# because the two coins swapped ranking positions in this run.
#
# Column D ("Price") is stored as plain text in this sheet (values like
# "43.642.10" use dots as thousands separators, not decimal points), so a
# replacement value that would otherwise parse as a genuine number (e.g.
# "1.01") is entered with a leading apostrophe - exactly like typing
# '1.01 into the cell in the Excel UI - so Excel keeps it as text instead
# of silently turning the cell into a numeric one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '43.642.10'; E = '  -0.09%  ' },
    @{ Row = 3; D = '2.279.77'; E = '  -0.33%  ' },
    @{ Row = 4; D = '1.01'; E = '  +0.38%  ' },
    @{ Row = 5; D = '112.71'; E = '  +9.73%  ' },
    @{ Row = 6; D = '266.94'; E = '  -1.41%  ' },
    @{ Row = 7; E = '  +0.28%  ' },
    @{ Row = 8; E = '  +0.11%  ' },
    @{ Row = 9; D = '0.612'; E = '  +1.19%  ' },
    @{ Row = 10; D = '48.36'; E = '  +5.22%  ' },
    @{ Row = 11; D = '0.0938'; E = '  +0.97%  ' },
    @{ Row = 12; E = '  +9.54%  ' },
    @{ Row = 13; D = '0.108'; E = '  +0.76%  ' },
    @{ Row = 14; D = '15.69'; E = '  +1.03%  ' },
    @{ Row = 15; D = '2.620.21'; E = '  -0.53%  ' },
    @{ Row = 16; D = '0.873'; E = '  +1.97%  ' },
    @{ Row = 17; D = '2.272.05'; E = '  -0.92%  ' },
    @{ Row = 18; D = '43.496.38'; E = '  -0.46%  ' },
    @{ Row = 19; D = '0.0000109'; E = '  -1.32%  ' },
    @{ Row = 20; D = '6.99'; E = '  +11.90%  ' },
    @{ Row = 21; D = '72.01' },
    @{ Row = 22; E = '  -3.52%  ' },
    @{ Row = 23; D = '9.94'; E = '  +7.71%  ' },
    @{ Row = 24; D = '232.23'; E = '  -0.41%  ' },
    @{ Row = 25; E = '  -0.87%  ' },
    @{ Row = 26; E = '  -0.01%  ' },
    @{ Row = 27; E = '  +2.80%  ' },
    @{ Row = 28; D = '41.38'; E = '  -0.24%  ' },
    @{ Row = 29; D = '3.39'; E = '  -1.64%  ' },
    @{ Row = 30; E = '  -1.49%  ' },
    @{ Row = 31; D = '173.38'; E = '  -2.23%  ' },
    @{ Row = 32; D = '21.52'; E = '  -1.17%  ' },
    @{ Row = 33; D = '0.0917'; E = '  +1.88%  ' },
    @{ Row = 34; E = '  +3.53%  ' },
    @{ Row = 35; E = '  +0.16%  ' },
    @{ Row = 36; D = '4.63'; E = '  -5.08%  ' },
    @{ Row = 37; E = '  -2.19%  ' },
    @{ Row = 38; E = '  -3.52%  ' },
    @{ Row = 39; D = '3.74'; E = '  +4.99%  ' },
    @{ Row = 40; D = '14.77'; E = '  +20.60%  ' },
    @{ Row = 41; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '2.43'; E = '  +4.64%  ' },
    @{ Row = 42; B = 'MultiversX'; C = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'; D = '74.30'; E = '  +13.05%  ' },
    @{ Row = 43; E = '  +0.20%  ' },
    @{ Row = 44; E = '  +19.35%  ' },
    @{ Row = 45; E = '  +0.15%  ' },
    @{ Row = 46; E = '  +0.67%  ' },
    @{ Row = 47; D = '8.67'; E = '  -1.29%  ' },
    @{ Row = 48; D = '102.02'; E = '  +3.33%  ' },
    @{ Row = 49; E = '  -1.90%  ' },
    @{ Row = 50; E = '  +1.75%  ' },
    @{ Row = 51; E = '  +1.97%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in 'B', 'C', 'D', 'E') {
        if (-not $u.ContainsKey($col)) { continue }
        $value = $u[$col]

        # Force plain numeric-looking Price cells to stay text.
        if ($col -eq 'D' -and $value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
            $value = "'" + $value
        }

        $ws.Range("$col$row").Value = $value
    }
}
